# Apply updated cryptos data (prices/volumes/row content) as of Sat Apr 13 16:48:53 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.453.17"
$ws.Range("E2").Value = "'  -1.78%  "
$ws.Range("D3").Value = "'3.251.65"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'587.31"
$ws.Range("E5").Value = "'  -4.44%  "
$ws.Range("D6").Value = "'146.57"
$ws.Range("E6").Value = "'  -11.90%  "
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D8").Value = "'3.241.61"
$ws.Range("E8").Value = "'  -5.09%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "'  -9.60%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "'  -13.76%  "
$ws.Range("D11").Value = "'6.66"
$ws.Range("E11").Value = "'  -4.22%  "
$ws.Range("D12").Value = "'0.495"
$ws.Range("E12").Value = "'  -11.49%  "
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("E13").Value = "'  -9.48%  "
$ws.Range("D14").Value = "'37.47"
$ws.Range("E14").Value = "'  -14.46%  "
$ws.Range("D15").Value = "'3.774.48"
$ws.Range("E15").Value = "'  -5.06%  "
$ws.Range("D16").Value = "'67.452.07"
$ws.Range("E16").Value = "'  -1.92%  "
$ws.Range("D17").Value = "'3.257.67"
$ws.Range("E17").Value = "'  -4.97%  "
$ws.Range("E18").Value = "'  -6.13%  "
$ws.Range("D19").Value = "'514.61"
$ws.Range("E19").Value = "'  -10.38%  "
$ws.Range("D20").Value = "'6.96"
$ws.Range("E20").Value = "'  -13.81%  "
$ws.Range("D21").Value = "'14.66"
$ws.Range("E21").Value = "'  -13.71%  "
$ws.Range("D22").Value = "'0.741"
$ws.Range("E22").Value = "'  -11.83%  "
$ws.Range("D23").Value = "'7.63"
$ws.Range("E23").Value = "'  -14.56%  "
$ws.Range("D24").Value = "'84.56"
$ws.Range("E24").Value = "'  -11.31%  "
$ws.Range("D25").Value = "'13.14"
$ws.Range("E25").Value = "'  -11.93%  "
$ws.Range("E26").Value = "'  +0.06%  "
$ws.Range("E27").Value = "'  -12.14%  "
$ws.Range("D28").Value = "'2.10"
$ws.Range("E28").Value = "'  -12.17%  "
$ws.Range("D29").Value = "'7.84"
$ws.Range("E29").Value = "'  -7.86%  "
$ws.Range("D30").Value = "'28.46"
$ws.Range("E30").Value = "'  -12.35%  "
$ws.Range("E31").Value = "'  -4.14%  "
$ws.Range("D32").Value = "'2.60"
$ws.Range("E32").Value = "'  -5.76%  "
$ws.Range("D33").Value = "'6.40"
$ws.Range("E33").Value = "'  -17.69%  "
$ws.Range("B34").Value = "'FirstDigitalUSD"
$ws.Range("C34").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  -0.03%  "
$ws.Range("B35").Value = "'NEARProtocol"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'5.56"
$ws.Range("E35").Value = "'  -14.78%  "
$ws.Range("D36").Value = "'55.99"
$ws.Range("E36").Value = "'  +0.03%  "
$ws.Range("D37").Value = "'502.71"
$ws.Range("E37").Value = "'  -15.05%  "
$ws.Range("D38").Value = "'0.0431"
$ws.Range("E38").Value = "'  -6.76%  "
$ws.Range("D39").Value = "'0.0836"
$ws.Range("E39").Value = "'  -11.91%  "
$ws.Range("B40").Value = "'Kaspa"
$ws.Range("C40").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.124"
$ws.Range("E40").Value = "'  -11.39%  "
$ws.Range("B41").Value = "'Cosmos"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.75"
$ws.Range("E41").Value = "'  -16.01%  "
$ws.Range("D42").Value = "'2.70"
$ws.Range("E42").Value = "'  -10.89%  "
$ws.Range("D43").Value = "'2.908.62"
$ws.Range("E43").Value = "'  -9.42%  "
$ws.Range("D44").Value = "'0.260"
$ws.Range("E44").Value = "'  -10.81%  "
$ws.Range("B45").Value = "'USDe"
$ws.Range("C45").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  -0.08%  "
$ws.Range("B46").Value = "'Fetch.AI"
$ws.Range("C46").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.14"
$ws.Range("E46").Value = "'  -9.51%  "
$ws.Range("B47").Value = "'InjectiveProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'26.06"
$ws.Range("E47").Value = "'  -15.86%  "
$ws.Range("D48").Value = "'0.0₃0556"
$ws.Range("E48").Value = "'  -17.30%  "
$ws.Range("D49").Value = "'124.28"
$ws.Range("E49").Value = "'  -6.28%  "
$ws.Range("D50").Value = "'0.111"
$ws.Range("E50").Value = "'  -10.89%  "
$ws.Range("D51").Value = "'2.25"
$ws.Range("E51").Value = "'  -18.51%  "

Write-Host "Applied 108 cell updates"
